# Update cryptos list with latest prices/volume (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '37.018.13'
$ws.Range("E2").Value = '  -0.26%  '

# Row 3
$ws.Range("D3").Value = '2.059.28'
$ws.Range("E3").Value = '  -2.14%  '

# Row 4
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.07%  '

# Row 5
$ws.Range("D5").Value = '249.03'
$ws.Range("E5").Value = '  -1.57%  '

# Row 6
$ws.Range("D6").Value = '0.670'
$ws.Range("E6").Value = '  +1.40%  '

# Row 7
$ws.Range("E7").Value = '  +0.06%  '

# Row 8
$ws.Range("D8").Value = '55.32'
$ws.Range("E8").Value = '  +12.12%  '

# Row 9
$ws.Range("D9").Value = '60.50'
$ws.Range("E9").Value = '  -0.63%  '

# Row 10
$ws.Range("D10").Value = '0.380'
$ws.Range("E10").Value = '  +0.78%  '

# Row 11
$ws.Range("E11").Value = '  +7.38%  '

# Row 12
$ws.Range("E12").Value = '  -0.96%  '

# Row 13
$ws.Range("E13").Value = '  +1.18%  '

# Row 14
$ws.Range("D14").Value = '2.360.81'
$ws.Range("E14").Value = '  -1.75%  '

# Row 15
$ws.Range("D15").Value = '0.813'
$ws.Range("E15").Value = '  -3.33%  '

# Row 16
$ws.Range("E16").Value = '  +1.44%  '

# Row 17
$ws.Range("D17").Value = '2.057.43'
$ws.Range("E17").Value = '  -4.45%  '

# Row 18
$ws.Range("D18").Value = '36.900.01'
$ws.Range("E18").Value = '  +0.16%  '

# Row 19
$ws.Range("B19").Value = 'Litecoin'
$ws.Range("C19").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D19").Value = '73.93'
$ws.Range("E19").Value = '  +0.70%  '

# Row 20
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").Value = '0.0₃0922'
$ws.Range("E20").Value = '  +11.73%  '

# Row 21
$ws.Range("D21").Value = '14.20'
$ws.Range("E21").Value = '  +6.40%  '

# Row 22
$ws.Range("D22").Value = '5.35'
$ws.Range("E22").Value = '  +0.00%  '

# Row 23
$ws.Range("D23").Value = '237.23'
$ws.Range("E23").Value = '  -1.96%  '

# Row 24
$ws.Range("E24").Value = '  -0.04%  '

# Row 25
$ws.Range("E25").Value = '  -4.05%  '

# Row 26
$ws.Range("D26").Value = '171.31'
$ws.Range("E26").Value = '  +0.34%  '

# Row 27
$ws.Range("D27").Value = '9.07'
$ws.Range("E27").Value = '  -4.62%  '

# Row 28
$ws.Range("D28").Value = '20.05'
$ws.Range("E28").Value = '  -5.33%  '

# Row 29
$ws.Range("D29").Value = '1.99'
$ws.Range("E29").Value = '  -0.66%  '

# Row 30
$ws.Range("E30").Value = '  +1.23%  '

# Row 31
$ws.Range("D31").Value = '4.58'
$ws.Range("E31").Value = '  +1.09%  '

# Row 32
$ws.Range("E32").Value = '  +0.83%  '

# Row 33
$ws.Range("D33").Value = '0.0627'
$ws.Range("E33").Value = '  +1.73%  '

# Row 34
$ws.Range("D34").Value = '4.36'
$ws.Range("E34").Value = '  +6.00%  '

# Row 35
$ws.Range("E35").Value = '  +0.14%  '

# Row 36
$ws.Range("D36").Value = '0.0875'
$ws.Range("E36").Value = '  -5.69%  '

# Row 37
$ws.Range("E37").Value = '  -5.76%  '

# Row 38
$ws.Range("E38").Value = '  +0.06%  '

# Row 39
$ws.Range("E39").Value = '  +0.01%  '

# Row 40
$ws.Range("D40").Value = '0.105'
$ws.Range("E40").Value = '  +23.11%  '

# Row 41
$ws.Range("E41").Value = '  +6.81%  '

# Row 42
$ws.Range("B42").Value = 'FTXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D42").Value = '4.60'
$ws.Range("E42").Value = '  +57.01%  '

# Row 43
$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").Value = '0.0225'
$ws.Range("E43").Value = '  -0.58%  '

# Row 44
$ws.Range("B44").Value = 'ARBITRUM'
$ws.Range("C44").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D44").Value = '1.15'
$ws.Range("E44").Value = '  -2.51%  '

# Row 45
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").Value = '96.76'
$ws.Range("E45").Value = '  -1.76%  '

# Row 46
$ws.Range("E46").Value = '  +0.59%  '

# Row 47
$ws.Range("E47").Value = '  +11.12%  '

# Row 48
$ws.Range("E48").Value = '  +6.68%  '

# Row 49
$ws.Range("D49").Value = '1.298.66'
$ws.Range("E49").Value = '  -3.65%  '

# Row 50
$ws.Range("E50").Value = '  -0.49%  '

# Row 51
$ws.Range("B51").Value = 'FraxShare'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D51").Value = '6.87'
$ws.Range("E51").Value = '  -3.42%  '
